$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Data edits: mask parent_id for NLO (row4) and SLO (row5) to 1 ---
$ws.Range("C4").Value = 1
$ws.Range("C5").Value = 1

# --- Clear the helper formula cells F75:G79 (CONCAT UNIT-A / UNIT-B) ---
$ws.Range("F75:G79").ClearContents()

# --- View state: scroll position + selection ---
$ws.Range("H54").Select()
$excel.ActiveWindow.ScrollRow = 37
